$wb = $excel.ActiveWorkbook

# Activate the StreamData sheet (this updates workbookView.activeTab
# and moves tabSelected from the previously active sheet to this one).
$ws = $wb.Worksheets.Item("StreamData")
$ws.Activate()

# Update the ExtremT (column L) values for rows 2-13.
$ws.Range("L2").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("L4").Value = 200
$ws.Range("L5").Value = 300
$ws.Range("L6").Value = 0
$ws.Range("L7").Value = 100
$ws.Range("L8").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("L10").Value = 400
$ws.Range("L11").Value = 200
$ws.Range("L12").Value = 0
$ws.Range("L13").Value = 100

# Move the selection within the frozen pane to M21.
$ws.Range("M21").Select()
